$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing rows (2..150) down to (3..151).
$ws.Rows("2:2").Insert()

# Populate the new row 2 with the new customer's data (a new discount/debt entry
# added to the top of the "Danh sách khách hàng" list).
$ws.Range("A2").Value = "KH"
$ws.Range("B2").Value = 376
$ws.Range("C2").Value = "Trần Thị Lệ"
$ws.Range("D2").Value = "CẦN THƠ"

# Phone numbers are stored as text (leading zero must be preserved) - force the
# cell to text format before assigning, then clear the number-format override so
# the cell keeps using the default style like the rest of the sheet.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "0919972597"
$ws.Range("F2").ClearFormats()

$ws.Range("I2").Value = 19000000
$ws.Range("J2").Value = 0
